$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.618.48'
$ws.Range("E2").Value = '  -1.44%  '
$ws.Range("D3").Value = '1.591.40'
$ws.Range("E3").Value = '  -1.76%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("E6").Value = '  -0.97%  '
$ws.Range("E7").Value = '  +0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.0617'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.51%  '
$ws.Range("E9").Value = '  -2.51%  '
$ws.Range("E10").Value = '  -2.41%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0834'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.58%  '
$ws.Range("D12").Value = '1.814.74'
$ws.Range("E12").Value = '  -1.74%  '
$ws.Range("D13").Value = '1.595.35'
$ws.Range("E13").Value = '  -2.11%  '
$ws.Range("E14").Value = '  -2.72%  '
$ws.Range("E15").Value = '  -2.77%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.91'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.35%  '
$ws.Range("D17").Value = '26.589.64'
$ws.Range("E17").Value = '  -1.49%  '
$ws.Range("E18").Value = '  -2.39%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '208.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.77%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.43%  '
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.29'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.95%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.88'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.43%  '
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("E27").Value = '  -3.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.115'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.29'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.76%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0504'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.16'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.48%  '
$ws.Range("E32").Value = '  -3.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.660'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -8.95%  '
$ws.Range("E34").Value = '  -3.39%  '
$ws.Range("D35").Value = '1.295.66'
$ws.Range("E35").Value = '  -3.43%  '
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("E37").Value = '  -5.08%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0172'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.830'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.02%  '
$ws.Range("E40").Value = '  +0.11%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.791'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.13%  '
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("E43").Value = '  -1.95%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = '1.727.35'
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '89.04'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.69%  '
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.816'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.03%  '
$ws.Range("E49").Value = '  -2.89%  '
$ws.Range("E50").Value = '  -2.42%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.53'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.67%  '

Write-Host "Applied all changes"
